# Update the "DateProd" timestamp column (B) on each sheet with the results
# of the latest Demo/QA CF verification run (per commit message:
# "Fixes for Demo Verification SCript and Demo/QA CF").

$wb = $excel.ActiveWorkbook

# CC-Payments-Auth sheet: rows 2-7
$wsAuth = $wb.Worksheets.Item("CC-Payments-Auth")
$wsAuth.Range("B2").Value = "Thu Jun 19 17:53:07 IST 2025"
$wsAuth.Range("B3").Value = "Thu Jun 19 17:53:47 IST 2025"
$wsAuth.Range("B4").Value = "Thu Jun 19 17:54:25 IST 2025"
$wsAuth.Range("B5").Value = "Thu Jun 19 17:55:05 IST 2025"
$wsAuth.Range("B6").Value = "Thu Jun 19 17:55:44 IST 2025"
$wsAuth.Range("B7").Value = "Thu Jun 19 17:56:21 IST 2025"

# ACH-Payments-Debit sheet: rows 2-10
$wsAch = $wb.Worksheets.Item("ACH-Payments-Debit")
$wsAch.Range("B2").Value = "Thu Jun 19 17:57:00 IST 2025"
$wsAch.Range("B3").Value = "Thu Jun 19 17:57:36 IST 2025"
$wsAch.Range("B4").Value = "Thu Jun 19 17:58:18 IST 2025"
$wsAch.Range("B5").Value = "Thu Jun 19 17:58:56 IST 2025"
$wsAch.Range("B6").Value = "Thu Jun 19 17:59:31 IST 2025"
$wsAch.Range("B7").Value = "Thu Jun 19 18:00:10 IST 2025"
$wsAch.Range("B8").Value = "Thu Jun 19 18:00:50 IST 2025"
$wsAch.Range("B9").Value = "Thu Jun 19 18:01:25 IST 2025"
$wsAch.Range("B10").Value = "Thu Jun 19 18:02:03 IST 2025"

# CC-Payments-Sale sheet: row 2
$wsSale = $wb.Worksheets.Item("CC-Payments-Sale")
$wsSale.Range("B2").Value = "Thu Jun 19 18:02:44 IST 2025"
